$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    3 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    4 = @(0.1169995834814548, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 0.8410281280571317)
    5 = @(3.272327238179451, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 4.260211312413533)
    6 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 5.582307763322248)
    7 = @(0.1169995834814548, 0.3048912486333797, 18.71679738969934, 0.5333859586016987, 19.67207418041587)
    8 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 8.656069925401464)
    9 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 24.14949828602258)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
